# Update VEDL financial statement figures for the latest reporting refresh.
# (Commit message: "Doing Updates for Financials")
#
# The workbook holds a pre-built Income Statement / Balance Sheet / Cash Flow
# Statement for VEDL. This pass refreshes the numeric figures across the
# seven reporting periods (columns D:J) to the newly restated source data,
# leaving all labels, formatting and layout untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 13212400
$ws.Range("E8").Value = 10370800
$ws.Range("F8").Value = 9247100
$ws.Range("G8").Value = 10607600
$ws.Range("H8").Value = 10487000
$ws.Range("I8").Value = 10444500
$ws.Range("J8").Value = 8648800
$ws.Range("D9").Value = 9484400
$ws.Range("E9").Value = 7721500
$ws.Range("F9").Value = 12663400
$ws.Range("G9").Value = 14401700
$ws.Range("H9").Value = 8067200
$ws.Range("I9").Value = 8049300
$ws.Range("J9").Value = 6304500
$ws.Range("D10").Value = 3728000
$ws.Range("E10").Value = 2649300
$ws.Range("F10").Value = -3416400
$ws.Range("G10").Value = -3794100
$ws.Range("H10").Value = 2419800
$ws.Range("I10").Value = 2395200
$ws.Range("J10").Value = 2344300
$ws.Range("D14").Value = -4800
$ws.Range("D17").Value = 9975300
$ws.Range("E17").Value = 8162100
$ws.Range("F17").Value = 13134200
$ws.Range("G17").Value = 14934200
$ws.Range("H17").Value = 8643000
$ws.Range("I17").Value = 8571800
$ws.Range("J17").Value = 7093900
$ws.Range("D18").Value = 3237100
$ws.Range("E18").Value = 2208700
$ws.Range("F18").Value = -3887200
$ws.Range("G18").Value = -4326600
$ws.Range("H18").Value = 1844100
$ws.Range("I18").Value = 1872700
$ws.Range("J18").Value = 1554800
$ws.Range("D20").Value = -859100
$ws.Range("E20").Value = 562900
$ws.Range("F20").Value = 553300
$ws.Range("G20").Value = 699200
$ws.Range("H20").Value = 464900
$ws.Range("I20").Value = 521400
$ws.Range("J20").Value = 242800
$ws.Range("D21").Value = 3461500
$ws.Range("E21").Value = 3661100
$ws.Range("F21").Value = -2127900
$ws.Range("G21").Value = -2020000
$ws.Range("H21").Value = 4072600
$ws.Range("I21").Value = 4088600
$ws.Range("J21").Value = 2681800
$ws.Range("D22").Value = 739200
$ws.Range("E22").Value = 796700
$ws.Range("F22").Value = 778600
$ws.Range("G22").Value = 876200
$ws.Range("H22").Value = 908100
$ws.Range("I22").Value = 807500
$ws.Range("J22").Value = 507900
$ws.Range("D23").Value = 1638800
$ws.Range("E23").Value = 1974800
$ws.Range("F23").Value = -4112500
$ws.Range("G23").Value = -4503700
$ws.Range("H23").Value = 1400800
$ws.Range("I23").Value = 1586600
$ws.Range("J23").Value = 1289700
$ws.Range("D24").Value = 960500
$ws.Range("E24").Value = 549900
$ws.Range("F24").Value = -1490200
$ws.Range("G24").Value = -1566300
$ws.Range("H24").Value = 501000
$ws.Range("I24").Value = -108500
$ws.Range("J24").Value = 111500
$ws.Range("D26").Value = 678300
$ws.Range("E26").Value = 1425000
$ws.Range("F26").Value = -2622300
$ws.Range("G26").Value = -2937400
$ws.Range("H26").Value = 899800
$ws.Range("I26").Value = 1695100
$ws.Range("J26").Value = 1178200
$ws.Range("D27").Value = 198200
$ws.Range("E27").Value = 795800
$ws.Range("F27").Value = -1809700
$ws.Range("G27").Value = -1855900
$ws.Range("H27").Value = 223600
$ws.Range("I27").Value = 901800
$ws.Range("J27").Value = 749200
$ws.Range("D32").Value = 859100
$ws.Range("E32").Value = -562900
$ws.Range("F32").Value = -553300
$ws.Range("G32").Value = -699200
$ws.Range("H32").Value = -464900
$ws.Range("I32").Value = -521400
$ws.Range("J32").Value = -242800
$ws.Range("D33").Value = 198200
$ws.Range("E33").Value = 795800
$ws.Range("F33").Value = -1809700
$ws.Range("G33").Value = -1855900
$ws.Range("H33").Value = 223600
$ws.Range("I33").Value = 901800
$ws.Range("J33").Value = 749200
$ws.Range("D35").Value = 198200
$ws.Range("E35").Value = 795800
$ws.Range("F35").Value = -1809700
$ws.Range("G35").Value = -1855900
$ws.Range("H35").Value = 223600
$ws.Range("I35").Value = 901800
$ws.Range("J35").Value = 749200
$ws.Range("D41").Value = 944700
$ws.Range("E41").Value = 2044900
$ws.Range("F41").Value = 600900
$ws.Range("G41").Value = 812000
$ws.Range("H41").Value = 1757200
$ws.Range("I41").Value = 2310700
$ws.Range("J41").Value = 1141500
$ws.Range("D42").Value = 4234700
$ws.Range("E42").Value = 6947500
$ws.Range("F42").Value = 7888000
$ws.Range("G42").Value = 6491700
$ws.Range("H42").Value = 5920700
$ws.Range("I42").Value = 3835300
$ws.Range("J42").Value = 2235100
$ws.Range("D43").Value = 863700
$ws.Range("E43").Value = 657500
$ws.Range("F43").Value = 1010300
$ws.Range("G43").Value = 1318000
$ws.Range("H43").Value = 1273600
$ws.Range("I43").Value = 2692100
$ws.Range("J43").Value = 1616900
$ws.Range("D44").Value = 1733500
$ws.Range("E44").Value = 1406500
$ws.Range("F44").Value = 1175000
$ws.Range("G44").Value = 1268400
$ws.Range("H44").Value = 1341700
$ws.Range("I44").Value = 1385700
$ws.Range("J44").Value = 661000
$ws.Range("D45").Value = 414800
$ws.Range("E45").Value = 387200
$ws.Range("F45").Value = 242800
$ws.Range("G45").Value = 260000
$ws.Range("H45").Value = 271300
$ws.Range("I45").Value = 194500
$ws.Range("J45").Value = 97700
$ws.Range("D46").Value = 8191400
$ws.Range("E46").Value = 11443700
$ws.Range("F46").Value = 10917000
$ws.Range("G46").Value = 10150200
$ws.Range("H46").Value = 10564500
$ws.Range("I46").Value = 9450900
$ws.Range("J46").Value = 5752300
$ws.Range("D47").Value = 927900
$ws.Range("E47").Value = 668500
$ws.Range("F47").Value = 467300
$ws.Range("G47").Value = 464500
$ws.Range("H47").Value = 463500
$ws.Range("I47").Value = 451100
$ws.Range("J47").Value = 591800
$ws.Range("D48").Value = 15152600
$ws.Range("E48").Value = 14112700
$ws.Range("F48").Value = 14177700
$ws.Range("G48").Value = 19297000
$ws.Range("H48").Value = 25045900
$ws.Range("I48").Value = 50848800
$ws.Range("J48").Value = 5487800
$ws.Range("D49").Value = 120100
$ws.Range("E49").Value = 92200
$ws.Range("F49").Value = 91100
$ws.Range("G49").Value = 95300
$ws.Range("H49").Value = 101000
$ws.Range("I49").Value = 348900
$ws.Range("J49").Value = 183000
$ws.Range("D52").Value = 1121600
$ws.Range("E52").Value = 1429700
$ws.Range("F52").Value = 1376300
$ws.Range("G52").Value = 1251200
$ws.Range("H52").Value = 1159900
$ws.Range("I52").Value = 741600
$ws.Range("J52").Value = 173200
$ws.Range("D54").Value = 25513600
$ws.Range("E54").Value = 27746700
$ws.Range("F54").Value = 27029400
$ws.Range("G54").Value = 31258200
$ws.Range("H54").Value = 37334800
$ws.Range("I54").Value = 34912000
$ws.Range("J54").Value = 12188200
$ws.Range("D57").Value = 1219600
$ws.Range("E57").Value = 905800
$ws.Range("F57").Value = 919400
$ws.Range("G57").Value = 848800
$ws.Range("H57").Value = 1215000
$ws.Range("I57").Value = 1037400
$ws.Range("J57").Value = 488000
$ws.Range("D58").Value = 4536100
$ws.Range("E58").Value = 5973800
$ws.Range("F58").Value = 2636500
$ws.Range("G58").Value = 2331400
$ws.Range("H58").Value = 2338600
$ws.Range("I58").Value = 2736400
$ws.Range("J58").Value = 780600
$ws.Range("D59").Value = 3948600
$ws.Range("E59").Value = 5037000
$ws.Range("F59").Value = 4277800
$ws.Range("G59").Value = 3088500
$ws.Range("H59").Value = 2944400
$ws.Range("I59").Value = 3636400
$ws.Range("J59").Value = 1017900
$ws.Range("D60").Value = 9704400
$ws.Range("E60").Value = 11916500
$ws.Range("F60").Value = 7833600
$ws.Range("G60").Value = 6268700
$ws.Range("H60").Value = 6497900
$ws.Range("I60").Value = 6130500
$ws.Range("J60").Value = 2286500
$ws.Range("D61").Value = 3873700
$ws.Range("E61").Value = 4810200
$ws.Range("F61").Value = 7140100
$ws.Range("G61").Value = 7488100
$ws.Range("H61").Value = 7915000
$ws.Range("I61").Value = 7563100
$ws.Range("J61").Value = 1077100
$ws.Range("D62").Value = 1048200
$ws.Range("E62").Value = 698700
$ws.Range("F62").Value = 831000
$ws.Range("G62").Value = 2627700
$ws.Range("H62").Value = 4668400
$ws.Range("I62").Value = 4298400
$ws.Range("J62").Value = 565200
$ws.Range("D66").Value = 16901800
$ws.Range("E66").Value = 19410000
$ws.Range("F66").Value = 21084100
$ws.Range("G66").Value = 23144500
$ws.Range("H66").Value = 27219100
$ws.Range("I66").Value = 25070400
$ws.Range("J66").Value = 5746600
$ws.Range("D72").Value = 4371000
$ws.Range("E72").Value = 5341400
$ws.Range("F72").Value = 2819400
$ws.Range("G72").Value = 4879500
$ws.Range("H72").Value = 6845800
$ws.Range("I72").Value = 6747800
$ws.Range("J72").Value = 3671700
$ws.Range("D76").Value = 8611900
$ws.Range("E76").Value = 8336700
$ws.Range("F76").Value = 5945400
$ws.Range("G76").Value = 8113800
$ws.Range("H76").Value = 10115800
$ws.Range("I76").Value = 9841600
$ws.Range("J76").Value = 6441600
$ws.Range("D81").Value = 198200
$ws.Range("E81").Value = 795800
$ws.Range("F81").Value = -1809700
$ws.Range("G81").Value = -1855900
$ws.Range("H81").Value = 223600
$ws.Range("I81").Value = 901800
$ws.Range("J81").Value = 749200
$ws.Range("D83").Value = 1082800
$ws.Range("E83").Value = 889000
$ws.Range("F83").Value = 1205100
$ws.Range("G83").Value = 1606400
$ws.Range("H83").Value = 1762500
$ws.Range("I83").Value = 1693300
$ws.Range("J83").Value = 883700
$ws.Range("D89").Value = 5038000
$ws.Range("E89").Value = 3067000
$ws.Range("F89").Value = 1552300
$ws.Range("G89").Value = 1817800
$ws.Range("H89").Value = 812600
$ws.Range("I89").Value = 1404200
$ws.Range("J89").Value = 2227800
$ws.Range("D91").Value = -1052000
$ws.Range("E91").Value = -777200
$ws.Range("F91").Value = -874400
$ws.Range("G91").Value = -1500100
$ws.Range("H91").Value = -1378200
$ws.Range("I91").Value = -1233700
$ws.Range("J91").Value = -1464000
$ws.Range("D94").Value = -754100
$ws.Range("E94").Value = -1155500
$ws.Range("F94").Value = -616200
$ws.Range("G94").Value = -635400
$ws.Range("H94").Value = -761000
$ws.Range("I94").Value = -2214900
$ws.Range("J94").Value = -7012200
$ws.Range("D96").Value = -2087200
$ws.Range("E96").Value = -75000
$ws.Range("F96").Value = -251000
$ws.Range("G96").Value = -150100
$ws.Range("H96").Value = -127000
$ws.Range("I96").Value = -138200
$ws.Range("J96").Value = -193500
$ws.Range("D100").Value = -5164300
$ws.Range("E100").Value = -710300
$ws.Range("F100").Value = -762200
$ws.Range("G100").Value = -1250000
$ws.Range("H100").Value = -90800
$ws.Range("I100").Value = 26800
$ws.Range("J100").Value = 5360400
$ws.Range("D101").Value = 12200
$ws.Range("E101").Value = -4300
$ws.Range("F101").Value = 3200
$ws.Range("G101").Value = 4900
$ws.Range("H101").Value = 6800
$ws.Range("I101").Value = 59900
$ws.Range("J101").Value = 15100
$ws.Range("D102").Value = -868300
$ws.Range("E102").Value = 1196900
$ws.Range("F102").Value = 177100
$ws.Range("G102").Value = -62700
$ws.Range("H102").Value = -32400
$ws.Range("I102").Value = -724000
$ws.Range("J102").Value = 591100
